$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the Price column as text first so purely-numeric-looking values
# (e.g. "1.00", "0.999") are preserved exactly instead of being
# auto-converted to numbers (which would drop trailing zeros).
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '58.367.40'
$ws.Range('E2').Value = '  -2.38%  '

$ws.Range('D3').Value = '2.576.50'
$ws.Range('E3').Value = '  -2.71%  '

$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').Value = '539.41'
$ws.Range('E5').Value = '  +0.50%  '

$ws.Range('D6').Value = '142.83'
$ws.Range('E6').Value = '  -1.66%  '

$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('D8').Value = '0.577'
$ws.Range('E8').Value = '  +0.89%  '

$ws.Range('D9').Value = '6.80'
$ws.Range('E9').Value = '  +1.24%  '

$ws.Range('E10').Value = '  -3.61%  '

$ws.Range('E11').Value = '  +2.76%  '

$ws.Range('D12').Value = '0.332'
$ws.Range('E12').Value = '  -1.99%  '

$ws.Range('D13').Value = '3.030.63'
$ws.Range('E13').Value = '  -2.90%  '

$ws.Range('D14').Value = '58.206.69'
$ws.Range('E14').Value = '  -2.48%  '

$ws.Range('D15').Value = '20.57'
$ws.Range('E15').Value = '  -3.10%  '

$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.0000132'
$ws.Range('E16').Value = '  -2.21%  '

$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.540.63'
$ws.Range('E17').Value = '  -3.90%  '

$ws.Range('E18').Value = '  +0.94%  '

$ws.Range('D19').Value = '334.81'
$ws.Range('E19').Value = '  -2.81%  '

$ws.Range('D20').Value = '10.02'
$ws.Range('E20').Value = '  -2.20%  '

$ws.Range('D21').Value = '6.12'
$ws.Range('E21').Value = '  -3.63%  '

$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.17%  '

$ws.Range('D23').Value = '66.44'
$ws.Range('E23').Value = '  -0.62%  '

$ws.Range('D24').Value = '0.419'
$ws.Range('E24').Value = '  +0.76%  '

$ws.Range('E25').Value = '  +0.08%  '

$ws.Range('E26').Value = '  -5.17%  '

$ws.Range('D27').Value = '7.01'
$ws.Range('E27').Value = '  -3.90%  '

$ws.Range('E28').Value = '  +0.04%  '

$ws.Range('D29').Value = '0.0₃0730'
$ws.Range('E29').Value = '  -2.51%  '

$ws.Range('D30').Value = '1.64'
$ws.Range('E30').Value = '  -1.22%  '

$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = '155.73'
$ws.Range('E31').Value = '  +3.64%  '

$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').Value = '5.91'
$ws.Range('E32').Value = '  +1.02%  '

$ws.Range('D33').Value = '18.85'
$ws.Range('E33').Value = '  -1.11%  '

$ws.Range('D34').Value = '3.87'
$ws.Range('E34').Value = '  -4.07%  '

$ws.Range('D35').Value = '36.91'
$ws.Range('E35').Value = '  -0.52%  '

$ws.Range('B36').Value = 'SuiNetwork'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D36').Value = '0.843'
$ws.Range('E36').Value = '  +2.70%  '

$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '1.09'
$ws.Range('E37').Value = '  -5.37%  '

$ws.Range('E38').Value = '  -2.96%  '

$ws.Range('E39').Value = '  -3.77%  '

$ws.Range('E40').Value = '  -0.60%  '

$ws.Range('D41').Value = '279.30'
$ws.Range('E41').Value = '  -5.66%  '

$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  +0.01%  '

$ws.Range('D43').Value = '0.588'
$ws.Range('E43').Value = '  -2.58%  '

$ws.Range('E44').Value = '  -0.87%  '

$ws.Range('D45').Value = '0.0529'
$ws.Range('E45').Value = '  -3.02%  '

$ws.Range('D46').Value = '0.0939'
$ws.Range('E46').Value = '  -1.70%  '

$ws.Range('D47').Value = '18.41'
$ws.Range('E47').Value = '  -4.92%  '

$ws.Range('D48').Value = '0.0226'
$ws.Range('E48').Value = '  -0.35%  '

$ws.Range('D49').Value = '1.908.12'
$ws.Range('E49').Value = '  -3.15%  '

$ws.Range('D50').Value = '17.79'
$ws.Range('E50').Value = '  -3.18%  '

$ws.Range('D51').Value = '4.37'
$ws.Range('E51').Value = '  -4.38%  '

# Restore the default (Normal) cell style on the Price column so the
# on-disk styling matches the original workbook (no stray number format).
$ws.Range('D2:D51').Style = 'Normal'
